$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Santa Barbara" label for row 6 (was blank)
$ws.Range("A6").Value = "Santa Barbara"

# Update the active selection to A7 (was D15)
$ws.Range("A7").Select()
